$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 6..17 down to 7..18 (bottom-up) so row 6 is free
# for the newly inserted claim record. Copying the used range (A:AA) keeps
# the row's existing cell styles intact instead of synthesizing new ones.
for ($r = 17; $r -ge 6; $r--) {
    $srcRange = $ws.Range("A" + $r + ":AA" + $r)
    $dstRange = $ws.Range("A" + ($r + 1) + ":AA" + ($r + 1))
    $srcRange.Copy($dstRange)
}

# Populate the newly freed row 6 with the new claim record.
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "ssurgwsoadev4-oci.opc.oracleoutsourcing.com"
$ws.Range("C6").Value = "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/cc/ClaimCenter.do"
$ws.Range("D6").Value = "rherner"
$ws.Range("E6").Value = "silverarrow"
$ws.Range("F6").Value = "04104018336"
$ws.Range("G6").Value = "Motor"
$ws.Range("H6").Value = "20/05/2021"
$ws.Range("I6").Value = "12:00"
$ws.Range("J6").Value = 32039201
$ws.Range("K6").Value = "aseguradosgw@gmail.com"
$ws.Range("L6").Value = "Daño Parcial"
$ws.Range("M6").Value = "Otros"
$ws.Range("N6").Value = "Calle"
$ws.Range("O6").Value = "CAPITAL"
$ws.Range("P6").Value = "CIUDAD AUTONOMA BUENOS AIRES"
$ws.Range("Q6").Value = "LIMA"
$ws.Range("R6").Value = 967
$ws.Range("S6").Value = "Pérdidas parciales"
$ws.Range("T6").Value = "No"
$ws.Range("U6").Value = "Parcial (Resto del Vehiculo)"
$ws.Range("V6").ClearContents()
$ws.Range("W6").Value = "Sí"
$ws.Range("X6").Value = "No"
$ws.Range("Y6").Value = "No"
$ws.Range("Z6").ClearContents()
$ws.Range("AA6").ClearContents()

# New claim record needs its own "Correo" hyperlink, just like the other
# rows covered by the aseguradosgw@gmail.com mailto link.
$ws.Hyperlinks.Add($ws.Range("K6"), "mailto:aseguradosgw@gmail.com", "", "", "aseguradosgw@gmail.com") | Out-Null

# Reflect the cell the user last landed on after the edit.
$ws.Application.ActiveWindow.ScrollColumn = 13
$ws.Range("AA7").Select() | Out-Null

Write-Output "done"
